# Atualização de bases das ligas, do dia: 07-03-2024 às 23:43
#
# Swap the data (all columns B..AC) between several row pairs that were
# re-ordered upstream, then append the newly scraped match as row 63.
# Column A (sequence id) and column E (match date) are intentionally left
# untouched - they stay tied to the row position, exactly like in the
# source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $range1 = "B$r1`:AC$r1"
    $range2 = "B$r2`:AC$r2"
    $v1 = $ws.Range($range1).Value2
    $v2 = $ws.Range($range2).Value2
    $ws.Range($range1).Value2 = $v2
    $ws.Range($range2).Value2 = $v1
}

Swap-Rows 2 3
Swap-Rows 4 5
Swap-Rows 12 13
Swap-Rows 46 47

# Append the new match scraped for this update as row 63.
# First clone the number formats (without values) from row 62 so the
# id column (A) and date column (E) reuse the existing styles instead of
# Excel minting brand-new ones.
$ws.Range("A62").Copy()
$ws.Range("A63").PasteSpecial(-4122)

$ws.Range("E62").Copy()
$ws.Range("E63").PasteSpecial(-4122)

$ws.Range("A63").Value2 = 61
$ws.Range("B63").Value2 = 7905567
$ws.Range("C63").Value2 = "Germany Landesliga"
$ws.Range("D63").Value2 = "Germany Landesliga"
$ws.Range("E63").Value2 = 45354.46875
$ws.Range("F63").Value2 = "SV Schlebusch"
$ws.Range("G63").Value2 = "SC Rheinbach"
$ws.Range("H63").Value2 = 0
$ws.Range("I63").Value2 = 0
$ws.Range("J63").Value2 = "D"
$ws.Range("K63").Value2 = 2.9
$ws.Range("L63").Value2 = 3.4
$ws.Range("M63").Value2 = 2.1
$ws.Range("N63").Value2 = 2.9
$ws.Range("O63").Value2 = 3.4
$ws.Range("P63").Value2 = 2.1
$ws.Range("Q63").Value2 = 0.25
$ws.Range("R63").Value2 = 1.9
$ws.Range("S63").Value2 = 1.9
$ws.Range("T63").Value2 = 3
$ws.Range("U63").Value2 = 1.85
$ws.Range("V63").Value2 = 1.95
$ws.Range("W63").Value2 = -1
$ws.Range("X63").Value2 = 2.4
$ws.Range("Y63").Value2 = -1
$ws.Range("Z63").Value2 = 0.45
$ws.Range("AA63").Value2 = -0.5
$ws.Range("AB63").Value2 = -1
$ws.Range("AC63").Value2 = 0.95
